# Eventuri "All Stars" import catalogue update
# - Rename the active "Global" worksheet tab to reflect the new export date
# - Rename header cells A1 (manufacturer -> Manufacturer) and L1 (Product -> Name)
# - Re-sync the Tags (P) and Keywords (Q) columns with the Name (L) column,
#   which restores the comma-separated tag text back to the plain product name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab title
$ws.Name = "Global - Sep-22-2023"

# Column headers
$ws.Range("A1").Value2 = "Manufacturer"
$ws.Range("L1").Value2 = "Name"

# Copy the Name column (L) into Tags (P) and Keywords (Q) for every data row
for ($r = 2; $r -le 156; $r++) {
    $name = $ws.Cells.Item($r, 12).Value2
    $ws.Cells.Item($r, 16).Value2 = $name
    $ws.Cells.Item($r, 17).Value2 = $name
}
